# Append 46 new transaction rows (rows 873-918) to the Association-Rules
# dataset sheet, matching the refreshed export produced by the (now
# Windows-path-fixed) Python data pipeline.
#
# Each entry is @(ngayxem date-text, masp, makh). Column A holds the
# "ngayxem" date written as plain text (shared string), exactly like the
# other ~870 existing rows above it — NOT an Excel date serial. We force
# that by briefly marking the cell as Text (NumberFormat "@") before the
# assignment, then clearing the format again so the cell keeps the
# worksheet's default (unstyled) appearance, just like every other data
# row already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2024-06-29", 50, 61),
    @("2024-06-29", 91, 61),
    @("2024-06-29", 56, 61),
    @("2024-06-29", 75, 62),
    @("2024-06-29", 137, 33),
    @("2024-06-29", 138, 33),
    @("2024-06-29", 136, 33),
    @("2024-06-29", 139, 33),
    @("2024-06-29", 140, 33),
    @("2024-06-29", 141, 33),
    @("2024-06-29", 142, 33),
    @("2024-06-29", 143, 33),
    @("2024-06-29", 144, 33),
    @("2024-06-29", 145, 33),
    @("2024-06-29", 146, 33),
    @("2024-06-29", 147, 33),
    @("2024-06-29", 138, 33),
    @("2024-06-29", 49, 33),
    @("2024-06-29", 59, 33),
    @("2024-06-29", 58, 33),
    @("2024-06-29", 80, 33),
    @("2024-06-29", 136, 33),
    @("2024-06-29", 75, 33),
    @("2024-06-29", 64, 33),
    @("2024-06-29", 138, 33),
    @("2024-06-29", 137, 67),
    @("2024-06-29", 182, 67),
    @("2024-06-29", 147, 67),
    @("2024-06-29", 167, 67),
    @("2024-06-29", 174, 67),
    @("2024-06-29", 226, 67),
    @("2024-06-29", 227, 45),
    @("2024-06-29", 58, 45),
    @("2024-06-29", 207, 45),
    @("2024-06-29", 139, 45),
    @("2024-06-29", 139, 45),
    @("2024-06-29", 55, 45),
    @("2024-06-29", 49, 45),
    @("2024-06-30", 212, 55),
    @("2024-06-30", 192, 27),
    @("2024-06-30", 49, 27),
    @("2024-06-30", 192, 27),
    @("2024-06-30", 62, 27),
    @("2024-07-01", 89, 68),
    @("2024-07-01", 228, 68),
    @("2024-07-01", 192, 52)
)

$startRow = 873
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
